$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New status entry for 5/2/2021 (row 27), following the pattern of the
# existing rows (e.g. row 26): date in col A, status text in col B,
# list of file names in col D.

$statusText = "1. I have completed pointers 5 programs and 2 files programs from given list of programs and pushed to GitHub`nhttps://github.com/gandepallipavani/C_Programs`n2. Gone through interview questions `n a. Size of structure without using sizeof operator`n b. If (0), if (-1)`n c. Call by value and call by reference`n3. Gone through testing ppt shared by Srinivasa`n4. I have completed 2 hacker rank programs today"
$filesText = "forinC.txt`nSumof5Digit.txt"

# Copy the formatting of row 26's cells into the new row 27 cells so the
# new row matches the existing look & feel (date format, wrap text, etc.)
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)

$ws.Range("B26").Copy()
$ws.Range("B27").PasteSpecial(-4122)

$ws.Range("D26").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Populate the new row's values
$ws.Range("A27").Value2 = 44318
$ws.Range("B27").Value2 = $statusText
$ws.Range("D27").Value2 = $filesText

# Match the row height used by similarly-sized entries (e.g. row 16/26)
$ws.Rows.Item(27).RowHeight = 225

# Update the view so the new row is visible/selected, like the author did
$ws.Range("E27").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
